# Export leaderboards as Excel
# Update the participants/leaderboard sheet with the latest standings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (top entry, no explicit rank) - totals/flag summary row
$ws.Range("B2").Value = "Ensan 3ayesh mn zaman"
$ws.Range("C2").Value = 620
$ws.Range("D2").Value = 1

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "moo moo"
$ws.Range("C3").Value = 150
$ws.Range("D3").Value = "Yes"

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Lolo Amr"
$ws.Range("C4").Value = 110
$ws.Range("D4").Value = "Yes"

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Alaa Amr"
$ws.Range("C5").Value = 50
$ws.Range("D5").Value = "Yes"

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Habiba Gamil"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = "Yes"

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Mariam Wael"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = "Yes"

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Salma Abosabie"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = "No"

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Britney Spears"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = "No"

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Mohamed Ahmed"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = "Yes"

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Mohamed Ahmed"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = "Yes"

# Row 12 - no Total points entry anymore
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Mohamed Ahmed"
$ws.Range("C12").ClearContents()
$ws.Range("D12").Value = "No"

# Row 13 - no Total points entry anymore
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Mohamed Ahmed"
$ws.Range("C13").ClearContents()
$ws.Range("D13").Value = "No"

# New row 14
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "Michael jackson"
$ws.Range("D14").Value = "No"
